$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 0.110552
$ws.Range("N2").Value = 0.331656
$ws.Range("O2").Value = 0.01126249561724847
$ws.Range("P2").Value = 0.01126249561724847
$ws.Range("Q2").Value = 0.146654082224
$ws.Range("R2").Value = 1.319886740016
$ws.Range("S2").Value = 0.01126249561724847
$ws.Range("T2").Value = 0.01126249561724847

$ws.Range("O3").Value = 0.9181055646724333
$ws.Range("P3").Value = 0.9181055646724334
$ws.Range("S3").Value = 0.9181055646724333
$ws.Range("T3").Value = 0.9181055646724334

$ws.Range("M4").Value = 0.6933189999999999
$ws.Range("N4").Value = 2.079957
$ws.Range("O4").Value = 0.07063193971031816
$ws.Range("P4").Value = 0.07063193971031817
$ws.Range("Q4").Value = 0.9197306392779999
$ws.Range("R4").Value = 8.277575753501999
$ws.Range("S4").Value = 0.07063193971031816
$ws.Range("T4").Value = 0.07063193971031817
